$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = "Proceso"
$ws.Range("O4").Value = "Afiliado"
$ws.Range("O12").Value = "Afiliado"
$ws.Range("O14").Value = "Proceso"
$ws.Range("O15").Value = "Afiliado"
$ws.Range("O16").Value = "Contacto"
$ws.Range("O21").Value = "Contacto"

$ws.Range("N6").Select()
